# Commit: "Replace $ in conditions to #"
# The IF/ELSE/ENDIF condition markers in cell D4 switch from the ${...}
# placeholder syntax to the #{...} syntax (the "price comparison" value
# stays the same, only the $ -> # prefix on IF/ELSE/ENDIF changes).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = '#{IF ${price} > 20}Expensive#{ELSE}Cheap#{ENDIF}'

# Selection moved from C5 to D5 in the saved view state.
$ws.Range("D5").Select()
